$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B:Q (cols 2..17) converge to the same metrics for every model row (2-26).
$rowValues = @(
    [double]"0.9994384710386303",
    [double]"0.9988584094093185",
    [double]"0.9999999000003338",
    [double]"0.9999985920531838",
    [double]"0.9999990556541269",
    [double]"0.0005241626371659963",
    [double]"0.001065624706365756",
    [double]"3.878798785116622e-08",
    [double]"1.090487133965852e-06",
    [double]"5.646375609085093e-07",
    [double]"0.0009694867605113856",
    [double]"0.02289459842770771",
    [double]"1.002695339014575",
    [double]"0.02386926888112749",
    [double]"73.107417091058",
    [double]"108.4548160122358"
)

# Row -> new A-column text (per the diff; unlisted rows keep their original label).
$labels = @{
    2 = "model_12_4_0"
    3 = "model_12_4_22"
    4 = "model_12_4_21"
    5 = "model_12_4_20"
    6 = "model_12_4_19"
    7 = "model_12_4_18"
    8 = "model_12_4_17"
    9 = "model_12_4_16"
    10 = "model_12_4_15"
    11 = "model_12_4_14"
    12 = "model_12_4_13"
    13 = "model_12_4_23"
    14 = "model_12_4_12"
    15 = "model_12_4_10"
    16 = "model_12_4_9"
    17 = "model_12_4_8"
    18 = "model_12_4_7"
    19 = "model_12_4_6"
    20 = "model_12_4_5"
    21 = "model_12_4_4"
    22 = "model_12_4_3"
    23 = "model_12_4_2"
    24 = "model_12_4_1"
    25 = "model_12_4_11"
    26 = "model_12_4_24"
}

foreach ($row in 2..26) {
    $ws.Cells.Item($row, 1).Value = $labels[$row]
    for ($i = 0; $i -lt $rowValues.Length; $i++) {
        $ws.Cells.Item($row, $i + 2).Value = $rowValues[$i]
    }
}
